{"js": "// The last paragraph of the document originally holds just an ellipsis\n// run (\"\u2026\") followed by a run with a single space (\" \"), then the\n// \"_GoBack\" bookmark. The edit extends that short placeholder paragraph\n// into the full meeting-notes paragraph. The \"_GoBack\" bookmark marks\n// where the author's typing left off, so it stays put right after the\n// first two (rewritten) runs and the rest of the text is appended after it.\n\nconst paragraphs = context.document.body.paragraphs;\nparagraphs.load(\"items\");\nawait context.sync();\n\nfor (let i = 0; i < paragraphs.items.length; i++) {\n  paragraphs.items[i].load(\"text\");\n}\nawait context.sync();\n\n// Locate the placeholder paragraph (its whole text is the ellipsis plus\n// a single space) instead of hard-coding an index.\nlet target = null;\nfor (let i = 0; i < paragraphs.items.length; i++) {\n  if (paragraphs.items[i].text === \"\\u2026 \") {\n    target = paragraphs.items[i];\n    break;\n  }\n}\nif (!target) {\n  throw new Error(\"Could not find the placeholder paragraph ('\u2026 ').\");\n}\n\n// Grab both existing runs *before* mutating anything, so the second\n// search is not confused by text inserted by the first edit.\nconst ellipsisResults = target.search(\"\\u2026\", { matchCase: true });\nconst spaceResults = target.search(\" \", { matchCase: true });\nellipsisResults.load(\"items\");\nspaceResults.load(\"items\");\nawait context.sync();\n\nconst ellipsisRun = ellipsisResults.items[0];\nconst spaceRun = spaceResults.items[0];\n\n// 1) Rewrite the two existing runs in place so their original run\n//    properties (w:lang = nl-NL) are preserved.\nellipsisRun.insertText(\n  \"Tijdens dit vrijwel korte gesprek hebben we het in\",\n  Word.InsertLocation.replace\n);\nawait context.sync();\n\nspaceRun.insertText(\n  \" eerste instantie over onze spoe\",\n  Word.InsertLocation.replace\n);\nawait context.sync();\n\n// 2) The \"_GoBack\" bookmark sits right after those two runs; insert the\n//    remainder of the text right after it, keeping it where Word left it.\nconst bookmark = context.document.getBookmarkRangeOrNullObject(\"_GoBack\");\nawait context.sync();\n\nconst remainder =\n  \"dcursus PHP -> LDAP \" +\n  \"gehad. Dit a\" +\n  \"angezien wij gisteren bij J. van der Veen zijn geweest, welke ons \" +\n  \"deze\" +\n  \" spoedc\" +\n  \"ursus PHP -> LDAP heeft gegeven.  Verder hebben we het over de server gehad, welke wij nog steeds geen informatie\" +\n  \" van\" +\n  \" hebben. Tot slot hebben we laten zien dat we onze site nu ook draaien hebben op onze webserv.nhl.nl website. Hiervoor heeft Yme na \" +\n  \"afloop van \" +\n  \"het gesprek bij J. Holwerda nog een database geregeld.\";\n\nbookmark.insertText(remainder, Word.InsertLocation.after);\nawait context.sync();\n", "ps1": "# The last paragraph of the document originally holds just an ellipsis\n# run (\"...\") followed by a run with a single space (\" \"), then the\n# \"_GoBack\" bookmark. The edit extends that short placeholder paragraph\n# into the full meeting-notes paragraph. The \"_GoBack\" bookmark marks\n# where the author's typing left off, so it stays put right after the\n# first two (rewritten) runs and the rest of the text is appended after it.\n\n$d = $word.ActiveDocument\n$ellipsis = [char]0x2026\n\n# Locate the placeholder paragraph (its whole text is the ellipsis plus a\n# single space, followed by the paragraph mark) instead of hard-coding an\n# index.\n$target = $null\nforeach ($p in $d.Paragraphs) {\n    if ($p.Range.Text -eq ($ellipsis + \" \" + [char]0x0D)) {\n        $target = $p\n        break\n    }\n}\nif ($target -eq $null) {\n    throw \"Could not find the placeholder paragraph ('... ').\"\n}\n\n$paraRange = $target.Range\n\n# Locate the ellipsis run within the paragraph.\n$ellipsisRange = $paraRange.Duplicate\n$ellipsisRange.Find.Execute($ellipsis)\n\n# 1) Rewrite the two existing runs in place so their original run\n#    properties (w:lang = nl-NL) are preserved. The space run is located\n#    *after* the ellipsis range has been rewritten, since this engine\n#    does not auto-shift independently created Range objects.\n$ellipsisRange.Text = \"Tijdens dit vrijwel korte gesprek hebben we het in\"\n\n$spaceRange = $d.Range($ellipsisRange.End, $ellipsisRange.End + 1)\n$spaceRange.Text = \" eerste instantie over onze spoe\"\n\n# 2) The \"_GoBack\" bookmark sits right after those two runs; insert the\n#    remainder of the text right after it, keeping it where Word left it.\n$bookmark = $d.Bookmarks.Item(\"_GoBack\")\n$remainder = \"dcursus PHP -> LDAP \" + `\n    \"gehad. Dit a\" + `\n    \"angezien wij gisteren bij J. van der Veen zijn geweest, welke ons \" + `\n    \"deze\" + `\n    \" spoedc\" + `\n    \"ursus PHP -> LDAP heeft gegeven.  Verder hebben we het over de server gehad, welke wij nog steeds geen informatie\" + `\n    \" van\" + `\n    \" hebben. Tot slot hebben we laten zien dat we onze site nu ook draaien hebben op onze webserv.nhl.nl website. Hiervoor heeft Yme na \" + `\n    \"afloop van \" + `\n    \"het gesprek bij J. Holwerda nog een database geregeld.\"\n\n$bookmark.Range.InsertAfter($remainder)\n"}
